$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Renewal / "Using Oldest Entry Date" rows -----------------------------
# Rows 3, 4 and 5 are new VIN-entry-date-symbol rows. They reuse the exact
# same layout/formatting as the existing row 2 (a SYMBOL_2000_CHOICE_T
# entry), so copy row 2 down first and then overwrite just the cells that
# differ for each new scenario.

$ws.Range("A2:AL2").Copy($ws.Range("A3:AL3"))
$ws.Range("A2:AL2").Copy($ws.Range("A4:AL4"))
$ws.Range("A2:AL2").Copy($ws.Range("A5:AL5"))

# Row 4's symbol columns use "C" - set this first so the new shared string
# table order matches (C before SYMBOL_2000_ENTRY_DATE).
$ws.Range("AE4:AH4").Value = "C"

# Row 3 - entry date 2000-01-01, not valid (oldest of the choice symbols)
$ws.Range("B3").Value = "SYMBOL_2000_ENTRY_DATE"
$ws.Range("AE3:AH3").Value = "N"
$ws.Range("AI3").Value = 20000101
$ws.Range("AJ3").Value = "N"

# Row 4 - entry date 2015-01-01
$ws.Range("B4").Value = "SYMBOL_2000_ENTRY_DATE"
$ws.Range("AI4").Value = 20150101
$ws.Range("AJ4").Value = "Y"

# Row 5 - entry date 2018-01-01
$ws.Range("B5").Value = "SYMBOL_2000_ENTRY_DATE"
$ws.Range("AE5:AH5").Value = "N"
$ws.Range("AI5").Value = 20180101
$ws.Range("AJ5").Value = "Y"

# --- Column B width now needs to fit the longer "SYMBOL_2000_ENTRY_DATE" text
$ws.Range("B1").EntireColumn.ColumnWidth = 24.1666666666667

# --- Selection / scroll position -------------------------------------------
[void]$ws.Range("B5").Select()
